# Applies the 05/10/2020 Philosophie course corrections:
#  1. Fix the spelling "cuasalité" -> "causalité" (and clear the stale
#     spell-check marker that surrounded it).
#  2. Merge the split "un"/"e" runs back into the word "une".
#  3. Remove the stale grammar-check marker around "En d'autre termes"
#     (sentence content is unchanged, only the proofing markup goes away).

$d = $word.ActiveDocument

# 1. "cuasalité" -> "causalité"
#    Match text spanning from the run before the typo through the run
#    after it so the stale <w:proofErr spellStart/spellEnd> wrapper
#    around the misspelled word is cleared together with the fix.
$d.Content.Find.Execute(
    "principe de cuasalité (vision rationnelle)", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "principe de causalité (vision rationnelle)", 2) | Out-Null

# 2. "Le déterminisme est un" + "e" + " doctrine..." -> "Le déterminisme est une doctrine..."
$d.Content.Find.Execute(
    "Le déterminisme est une doctrine philosophique et scientifique, le fatalisme est la croyance du destin, c’est-à-dire au fait qu’un événement est fixé par avance dans une logique de boules de villards ou de dominos.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Le déterminisme est une doctrine philosophique et scientifique, le fatalisme est la croyance du destin, c’est-à-dire au fait qu’un événement est fixé par avance dans une logique de boules de villards ou de dominos.",
    2) | Out-Null

# 3. Drop the gramStart/gramEnd proofing markers around "En d’autre termes"
#    by re-asserting the whole sentence (content itself is unchanged).
$d.Content.Find.Execute(
    "Inconvénient de cette démarche : elle peut sembler circulaire. Il faut déjà savoir que Socrate est mortel pour prétendre que tous les hommes sont mortels. En d’autre termes, la conclusion semble déjà contenue dans la première prémisse.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Inconvénient de cette démarche : elle peut sembler circulaire. Il faut déjà savoir que Socrate est mortel pour prétendre que tous les hommes sont mortels. En d’autre termes, la conclusion semble déjà contenue dans la première prémisse.",
    2) | Out-Null
